# "Command Line implementation added"
# - Features!B2 ("Run" flag for "Flight Search Feature") flips from Y to N.
# - Selection on Features sheet moves to B2.
# - Selection on Scenarios sheet moves to B5, and Scenarios becomes the active/visible tab.
# - Tags sheet loses the active-tab flag (it was previously selected).

$wb = $excel.ActiveWorkbook

$wsFeatures  = $wb.Worksheets.Item("Features")
$wsScenarios = $wb.Worksheets.Item("Scenarios")
$wsTags      = $wb.Worksheets.Item("Tags")

# Flip the Run flag for the second feature row.
$wsFeatures.Range("B2").Value = "N"

# Update selections on each sheet.
$wsFeatures.Range("B2").Select()
$wsScenarios.Range("B5").Select()
$wsTags.Range("D18").Select()

# Make Scenarios the active sheet/tab.
$wsScenarios.Activate()
